$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values to reflect the "Added Modified Reg iProctor P2,P3 TC's" change
$ws.Range("A2").Value = "XwjNd171"
$ws.Range("B2").Value = 23081822
$ws.Range("C2").Value = "usbmybz49"
$ws.Range("D2").Value = "vTd2`$D#6"
$ws.Range("F2").Value = "gBzTfrvq"
$ws.Range("G2").Value = "Gcvw"
